$d = $word.ActiveDocument

$d.Content.Find.Execute("Username: josnoble113@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Username: aobrien516@c2ken.net", 2)

$d.Content.Find.Execute("Password: jn11jn11", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Password: password1", 2)

$d.Content.Find.Execute("Firstname: josh", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Firstname: anthony", 2)

$d.Content.Find.Execute("Surname: noble", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Surname: obrien", 2)

$d.Content.Find.Execute("Address: 113 srtett", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Address: 113 road", 2)

$d.Content.Find.Execute("Postcode: de34 4ed", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Postcode: bt45 7yt", 2)

$d.Content.Find.Execute("Age: 19", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Age: 31", 2)

$d.Content.Find.Execute("Group: 4", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Group: 7", 2)
